# Apply corrections to official district/address names as per the commit:
# "corrected most names to the official names from website"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clear out the stray, empty Address (column F) cells that were left
#    behind with no content - rows for teachers whose address info was
#    already merged elsewhere. Clearing removes the now-superfluous empty
#    cells (and collapses the sheet's used range from A1:H47 to A1:G47).
$emptyFRows = @(4, 6, 7, 8, 9, 10, 14, 16, 18, 21, 22, 23, 25, 26, 29, 30, 33, 34, 37, 38, 41, 42, 43, 44, 45, 46, 47)
foreach ($r in $emptyFRows) {
    $ws.Cells.Item($r, 6).ClearContents()
}

# 2) Correct District (column G) names to their official spellings.
$ws.Range("G5").Value = "Tumakuru (Tumkur)"
$ws.Range("G17").Value = "Ballari (Bellary)"
$ws.Range("G32").Value = "Ballari (Bellary)"
$ws.Range("G39").Value = "Ballari (Bellary)"
$ws.Range("G40").Value = "Ballari (Bellary)"
